{"js": "// \"...en el sistema. Por otro lado, los procesos cooperativos...\" ->\n// \"...en el sistema, y tampoco puede comunicarse con otro proceso. Por otro lado, ...\"\n// Find the exact original phrase and replace it with the expanded sentence.\nconst body = context.document.body;\nconst results = body.search(\"sistema. Por otro lado,\", { matchCase: true });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Target phrase \"sistema. Por otro lado,\" not found.');\n}\n\nconst target = results.items[0];\ntarget.insertText(\n  \"sistema, y tampoco puede comunicarse con otro proceso. Por otro lado,\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# Split \"sistema. Por otro lado,\" into \"sistema\" + new sentence + \"Por otro lado,\"\n# by replacing the whole phrase with the expanded text.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"sistema. Por otro lado,\"\n$find.Replacement.Text = \"sistema, y tampoco puede comunicarse con otro proceso. Por otro lado,\"\n$find.Forward = $true\n$find.Wrap = 0            # wdFindStop\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdReplaceOne = 1 ; wdFindStop = 0\n$found = $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 0, $false, $find.Replacement.Text, 1)\n\nif (-not $found) {\n    throw 'Target phrase \"sistema. Por otro lado,\" not found.'\n}\n"}
